$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns B:I, rows 2-25 (contingencies with rene fine)
$dataBI = @(
  @(2.597803197262351, 0.7119485641435404, 0.03291285572721492, 1.361407321067915, 0.6798212621317248, 0.0007871586232472039, 0.01199277220940109, 0.005767783202876409),
  @(2.26066291172782, 0.6281112901472454, 0.02978638338481687, 1.18640414660436, 0.6162561042387509, 0.0007907025278933496, 0.008635475661823488, 0.003802278604460696),
  @(2.053033621404154, 0.5769465878343851, 0.02786294938359291, 1.079243740448391, 0.5777840511071304, 0.0007929509702298748, 0.006788786149038328, 0.002808140480463628),
  @(1.966867226263588, 0.5570060540242423, 0.02711732349066409, 1.03560817179681, 0.561636489855097, 0.0007938913891312019, 0.006084113324819063, 0.002507873075721534),
  @(1.950852273096388, 0.5547293813975784, 0.02704163866897247, 1.028341575679136, 0.5582333237960952, 0.0007940555458227334, 0.005967871774951172, 0.00253600607557658),
  @(2.047220876175572, 0.579499765070608, 0.02798496047145704, 1.078590135761104, 0.5755708414071776, 0.0007929820356354663, 0.006772998661915941, 0.002993613051261335),
  @(2.475449915545028, 0.6867448461070467, 0.03201250899027031, 1.300903148091805, 0.655126001160383, 0.0007883887738010078, 0.01077826887721478, 0.00526474178485703),
  @(3.31937571960276, 0.89486404983694, 0.03964780683899249, 1.74031629218878, 0.819821406831295, 0.0007799115522071909, 0.02055750911132859, 0.01132324856504763),
  @(3.917712953324951, 1.045199349042434, 0.04618919566138402, 1.963957886527425, 0.93232782730243, 0.0007741951350906243, 0.02843873462482227, 0.01710822360147279),
  @(4.026737869601959, 1.062850731216088, 0.05755877076850879, 1.273281716279527, 0.8817242950077855, 0.0007731507858902821, 0.04364912973931112, 0.01881957031504289),
  @(3.998808706553689, 1.04164719867569, 0.06634448787927028, 0.7757246460136571, 0.8193973257164942, 0.0007732720225684575, 0.07933126607416341, 0.01876606120865443),
  @(3.856906087453297, 0.9941221507705222, 0.07374033173330474, 0.3981850151505242, 0.7415224755416716, 0.0007743164300377484, 0.1320877733521399, 0.01762405812788259),
  @(3.705445098923747, 0.9502030679651909, 0.07837039433519521, 0.2067866589870491, 0.6801878403898911, 0.0007754502749329242, 0.1794594729324785, 0.01645342251417592),
  @(3.642840219142954, 0.9350401674155364, 0.07915881600541752, 0.1691515420980636, 0.6612674466339854, 0.0007759511119717162, 0.1913068405711869, 0.01602642478418481),
  @(3.415509897697746, 0.8831420113135096, 0.0746227652443352, 0.1649757372158618, 0.6290727610731253, 0.000778154564283884, 0.1761469891596477, 0.01389576357840028),
  @(3.321024665012828, 0.8665758524583111, 0.068686150943293, 0.2495924249046269, 0.6365869619460085, 0.0007792633737742622, 0.1374730480102357, 0.01290158642736738),
  @(3.333655881952268, 0.8764144046545539, 0.06103533650605186, 0.4836990774361354, 0.6812400819957105, 0.000779451795486122, 0.08540656275748404, 0.0125143417278295),
  @(3.423425903977261, 0.9129069093284556, 0.05334773736532838, 0.9192880719533036, 0.7509527051255276, 0.0007787894674839167, 0.04300193645821793, 0.01316157143171637),
  @(3.744961384763428, 1.014395815254744, 0.04493882850729136, 1.900556197826361, 0.8957133485068596, 0.0007757310825716802, 0.02618642432198026, 0.01606431092290084),
  @(4.234385973037206, 1.13884384872128, 0.04853039978854667, 2.232452845818756, 1.002870546799215, 0.0007711437311507585, 0.03398707486947439, 0.02118666147670556),
  @(4.551661220263384, 1.213745153735829, 0.05114237283794409, 2.398415709812426, 1.070752524721385, 0.000768271936691179, 0.03907437156944749, 0.02456729336757402),
  @(4.388043389984887, 1.170230267484214, 0.04958402825251795, 2.309816508383619, 1.03681460987049, 0.0007697827481974752, 0.03633694744138571, 0.02250664103468569),
  @(3.759689213248521, 1.012181686639536, 0.04394261903344443, 1.976461631646146, 0.9067763043870798, 0.0007756338485891714, 0.02669420700939273, 0.01574906392351672),
  @(3.082796153344248, 0.8434117984034515, 0.03782965901152835, 1.620805740526038, 0.7709755821200588, 0.000782176244044853, 0.01766107065060085, 0.009773803154761396)
)

$rows = $dataBI.Count
$cols = $dataBI[0].Count
$arrBI = New-Object 'object[,]' $rows,$cols
for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $arrBI[$r,$c] = $dataBI[$r][$c]
    }
}
$ws.Range("B2:I25").Value = $arrBI

# Updated values for columns P:Q, rows 2-25
$dataPQ = @(
  @(0.7161795053628026, 2.014744168414211),
  @(0.7252262043596858, 1.84802363948063),
  @(0.7314904468061201, 1.747272484144673),
  @(0.734578658242711, 1.704511633347835),
  @(0.7355487681603776, 1.694901505474832),
  @(0.7327682632237398, 1.739771223154349),
  @(0.7208571998837741, 1.947695453102057),
  @(0.7010876780096567, 2.382696616817753),
  @(0.696941688860889, 2.671041276445521),
  @(0.7518756670614692, 2.455349917182559),
  @(0.8018015336147926, 2.230018738977265),
  @(0.8516298066913066, 1.974162583525896),
  @(0.8871022934417567, 1.782298678331358),
  @(0.8957326648395139, 1.726423683432074),
  @(0.8846610237793584, 1.655776269964292),
  @(0.8577170160333907, 1.70188707370113),
  @(0.8143454827812491, 1.86245199742703),
  @(0.7672507512976097, 2.100504158686675),
  @(0.7024029010301973, 2.571054206805428),
  @(0.691096490122753, 2.862923816794051),
  @(0.6855652115969377, 3.045853247365073),
  @(0.6867202211229895, 2.95620379021554),
  @(0.6955573263841259, 2.610298883614547),
  @(0.7082226143536374, 2.250509085723536)
)

$rows2 = $dataPQ.Count
$cols2 = $dataPQ[0].Count
$arrPQ = New-Object 'object[,]' $rows2,$cols2
for ($r = 0; $r -lt $rows2; $r++) {
    for ($c = 0; $c -lt $cols2; $c++) {
        $arrPQ[$r,$c] = $dataPQ[$r][$c]
    }
}
$ws.Range("P2:Q25").Value = $arrPQ
